$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of sample data (row 4), mirroring the existing rows 2-3.
# Force text format first so the numeric-looking zip/state values
# ("222"/"2222") are stored as text rather than being coerced to numbers.
$ws.Range("A4:H4").NumberFormat = "@"

$ws.Range("A4").Value = "f3 address"
$ws.Range("B4").Value = " "
$ws.Range("C4").Value = "f2. first"
$ws.Range("D4").Value = "f2 last"
$ws.Range("E4").Value = " "
$ws.Range("F4").Value = "f2 city"
$ws.Range("G4").Value = "222"
$ws.Range("H4").Value = "2222"

# Reset to the default style so the new row matches the unstyled
# data rows above it (only the header row keeps the bold style).
$ws.Range("A4:H4").Style = "Normal"
